$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("tags")

# Add the new row 15, reusing the existing "dl" shared string so it
# stays referenced (it moves from row 13 down to row 15).
$ws.Range("A15").Value = "dl"
$ws.Range("B15").Value = "draw"
$ws.Range("C15").Value = 4

# Swap the "mtl" / "rmtl" rows (row 11 <-> row 12 values).
$ws.Range("A11").Value = "rmtl"
$ws.Range("A12").Value = "mtl"

# Row 13 now holds a brand-new tag "omtl" (filter for plots).
$ws.Range("A13").Value = "omtl"

# Match the saved selection (active cell A13) from the source workbook.
$ws.Activate() | Out-Null
$ws.Range("A13").Select() | Out-Null
